$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing team-name cells so re-assignment does not collide mid-way
$ws.Range("B2:B30").ClearContents()

$teams = @("POR", "NJN", "CLE", "DAL", "MIA", "SEA", "ATL", "WAS", "MIL", "LAC", "VAN", "DET", "SAS", "ORL", "UTA", "HOU", "DEN", "LAL", "GSW", "IND", "CHI", "PHI", "CHH", "BOS", "TOR", "SAC", "PHO", "NYK", "MIN")
for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $teams[$i]
}

# Update column C values (recomputed stats)
$ws.Cells.Item(2, 3).Value = 8.828571428571433
$ws.Cells.Item(3, 3).Value = 11.49230769230769
$ws.Cells.Item(4, 3).Value = 11.39333333333333
$ws.Cells.Item(5, 3).Value = 11.77857142857143
$ws.Cells.Item(6, 3).Value = 14.28571428571429
$ws.Cells.Item(7, 3).Value = 14.44285714285714
$ws.Cells.Item(8, 3).Value = 13.1375
$ws.Cells.Item(9, 3).Value = 11.9125
$ws.Cells.Item(10, 3).Value = 12.2875
$ws.Cells.Item(11, 3).Value = 13.35833333333333
$ws.Cells.Item(12, 3).Value = 13.82727272727273
$ws.Cells.Item(13, 3).Value = 13.51666666666667
$ws.Cells.Item(14, 3).Value = 12.91428571428571
$ws.Cells.Item(15, 3).Value = 13.11538461538461
$ws.Cells.Item(16, 3).Value = 15.2
$ws.Cells.Item(17, 3).Value = 14.12857142857143
$ws.Cells.Item(18, 3).Value = 12.42142857142857
$ws.Cells.Item(19, 3).Value = 15.76923076923077
$ws.Cells.Item(20, 3).Value = 10.36
$ws.Cells.Item(21, 3).Value = 11.41428571428571
$ws.Cells.Item(22, 3).Value = 13.73571428571429
$ws.Cells.Item(23, 3).Value = 15.3
$ws.Cells.Item(24, 3).Value = 11.85833333333333
$ws.Cells.Item(25, 3).Value = 12.775
$ws.Cells.Item(26, 3).Value = 10.98333333333333
$ws.Cells.Item(27, 3).Value = 11.84615384615385
$ws.Cells.Item(28, 3).Value = 15.075
$ws.Cells.Item(29, 3).Value = 13.30714285714286
$ws.Cells.Item(30, 3).Value = 12.92666666666667
